$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 385, pushing the existing row 385 (and everything
# below it) down by one. This grows the used range from A1:T468 to A1:T469.
$ws.Rows.Item(385).Insert()

# Populate the newly inserted row 385 with the new price-record data.
# Columns A, B, C, E, F, G, H, I, J carry the same market/category
# metadata as the row that used to sit at 385 (now row 386).
$ws.Cells.Item(385, 1).Value = 3
$ws.Cells.Item(385, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(385, 3).Value = "Coquimbo"
$ws.Cells.Item(385, 4).Value = 44508
$ws.Cells.Item(385, 5).Value = 5
$ws.Cells.Item(385, 6).Value = "Fruta"
$ws.Cells.Item(385, 7).Value = 100103
$ws.Cells.Item(385, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(385, 9).Value = 100103006
$ws.Cells.Item(385, 10).Value = "Nectarín"
$ws.Cells.Item(385, 11).Value = "Early Glo"
$ws.Cells.Item(385, 12).Value = "Primera"
$ws.Cells.Item(385, 13).Value = 60
$ws.Cells.Item(385, 14).Value = 22000
$ws.Cells.Item(385, 15).Value = 22000
$ws.Cells.Item(385, 16).Value = 22000
$ws.Cells.Item(385, 17).Value = "`$/caja 15 kilos granel"
$ws.Cells.Item(385, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(385, 19).Value = 1467
$ws.Cells.Item(385, 20).Value = 15
